$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns keep their text formatting so Excel does not
# auto-convert numeric-looking / percentage-looking strings.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.027.97"
$ws.Range("D3").Value = "1.897.59"
$ws.Range("E3").Value = "  +1.46%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "312.35"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "0.5023"
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").Value = "0.3918"
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("D9").Value = "0.09389"
$ws.Range("E9").Value = "  -2.46%  "
$ws.Range("D10").Value = "1.132"
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("D11").Value = "41.87"
$ws.Range("E11").Value = "  +2.25%  "
$ws.Range("D12").Value = "6.362"
$ws.Range("E12").Value = "  -1.53%  "
$ws.Range("D13").Value = "20.74"
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("D14").Value = "1.892.75"
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "7.298"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "0.00001116"
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("D18").Value = "92.47"
$ws.Range("E18").Value = "  -0.77%  "
$ws.Range("D19").Value = "0.06584"
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("D20").Value = "17.83"
$ws.Range("E20").Value = "  +2.05%  "
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "6.227"
$ws.Range("E22").Value = "  +1.09%  "
$ws.Range("D23").Value = "28.073.77"
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("D24").Value = "11.33"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "2.317"
$ws.Range("E25").Value = "  +1.69%  "
$ws.Range("D26").Value = "2.625"
$ws.Range("E26").Value = "  +2.65%  "
$ws.Range("D27").Value = "2.112.63"
$ws.Range("E27").Value = "  +1.30%  "
$ws.Range("D28").Value = "20.87"
$ws.Range("E28").Value = "  -1.25%  "
$ws.Range("D29").Value = "157.12"
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("D30").Value = "126.91"
$ws.Range("E30").Value = "  -0.40%  "
$ws.Range("D31").Value = "1.081"
$ws.Range("E31").Value = "  +2.01%  "
$ws.Range("D32").Value = "0.1064"
$ws.Range("E32").Value = "  +0.97%  "
$ws.Range("D33").Value = "5.620"
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("D34").Value = "3.619"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("D35").Value = "9.581"
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("D36").Value = "0.06615"
$ws.Range("E36").Value = "  -2.12%  "
$ws.Range("D37").Value = "0.02421"
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "0.2175"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "1.288"
$ws.Range("E39").Value = "  +9.50%  "
$ws.Range("E40").Value = "  -2.61%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.6394"
$ws.Range("E41").Value = "  +1.55%  "
$ws.Range("B42").Value = "InternetComputer(DFINITY)"
$ws.Range("C42").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D42").Value = "4.999"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").Value = "11.43"
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").Value = "13.35"
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("D46").Value = "0.5987"
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("D47").Value = "3.714"
$ws.Range("E47").Value = "  +1.58%  "
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("D49").Value = "2.025"
$ws.Range("E49").Value = "  +2.41%  "
$ws.Range("D50").Value = "123.06"
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("D51").Value = "1.176"
